$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.529578333333333
$ws.Range("H2").Value = 13.588735
$ws.Range("I2").Value = 0.2308013058217703
$ws.Range("J2").Value = 0.2308013058217703
$ws.Range("M2").Value = 0.467769
$ws.Range("Q2").Value = 2.118796327405
$ws.Range("R2").Value = 19.069166946645
$ws.Range("S2").Value = 0.2308013058217703
$ws.Range("T2").Value = 0.2308013058217703

# Row 3
$ws.Range("I3").Value = 0.3456045794970084
$ws.Range("J3").Value = 0.3456045794970085
$ws.Range("M3").Value = 0.467769
$ws.Range("Q3").Value = 3.172710445313
$ws.Range("S3").Value = 0.3456045794970084
$ws.Range("T3").Value = 0.3456045794970085

# Row 4
$ws.Range("G4").Value = 2.766332333333333
$ws.Range("H4").Value = 8.298997
$ws.Range("I4").Value = 0.1409564131327128
$ws.Range("J4").Value = 0.1409564131327128
$ws.Range("M4").Value = 0.467769
$ws.Range("Q4").Value = 1.294004509231
$ws.Range("R4").Value = 11.646040583079
$ws.Range("S4").Value = 0.1409564131327128
$ws.Range("T4").Value = 0.1409564131327128

# Row 5
$ws.Range("G5").Value = 5.546890666666666
$ws.Range("H5").Value = 16.640672
$ws.Range("I5").Value = 0.2826377015485084
$ws.Range("J5").Value = 0.2826377015485084
$ws.Range("M5").Value = 0.467769
$ws.Range("Q5").Value = 2.594663500256
$ws.Range("R5").Value = 23.351971502304
$ws.Range("S5").Value = 0.2826377015485084
$ws.Range("T5").Value = 0.2826377015485084
